# Weekly update: insert a new weekly record at row 173 for
# "Hortaliza, Macroferia Regional de Talca - Apio", pushing the
# existing historical rows (old 173-214) down to (174-215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 173; Excel shifts rows
# 173..214 down to 174..215 automatically (values, styles, formulas).
$ws.Rows.Item(173).Insert()

# Populate the newly-opened row 173 with this week's record. All
# columns other than D (Fecha), K/L/M (Precio min/max/promedio) and
# P (Precio $/Kg) are identical to the neighbouring rows for this
# market/category, so copy them across explicitly.
$ws.Cells.Item(173, 1).Value = 5
$ws.Cells.Item(173, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(173, 3).Value = "Maule"
$ws.Cells.Item(173, 4).Value = 44782
$ws.Cells.Item(173, 5).Value = 7
$ws.Cells.Item(173, 6).Value = 100112017
$ws.Cells.Item(173, 7).Value = "Apio"
$ws.Cells.Item(173, 8).Value = "Americana (o)"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 600
$ws.Cells.Item(173, 11).Value = 10000
$ws.Cells.Item(173, 12).Value = 10000
$ws.Cells.Item(173, 13).Value = 10000
$ws.Cells.Item(173, 14).Value = "$/docena de matas"
$ws.Cells.Item(173, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(173, 16).Value = 1667
$ws.Cells.Item(173, 17).Value = 6
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Make sure the number format used for the date column (D) matches
# the rest of the column.
$ws.Range("D173").NumberFormat = $ws.Range("D174").NumberFormat
